# Auto-generated edit script: rewrites the data tables on each sheet
# to match the COMM-revised column order, renamed headers, and renamed
# row labels described in the commit message / diff.

$wb = $excel.ActiveWorkbook

function Set-TextValue {
    param($Range, [string]$Text)
    # Force literal-text entry so numeric-looking strings (percentages,
    # dollar amounts, comma-separated counts, ...) are not auto-converted
    # into numbers/percentages by Excel.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# ---- Sheet: Overall ----
$ws = $wb.Worksheets.Item("Overall")
$ws.Range("A1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("B1").Value = "Number of 990 filers with government grants"
$ws.Range("C1").Value = "Total government grants (`$)"
$ws.Range("D1").Value = "Size of operating surplus with government grants"
$ws.Range("E1").Value = "Size of operating surplus without government grants"

Set-TextValue $ws.Range("A2") "66.00%"
Set-TextValue $ws.Range("B2") "547"
Set-TextValue $ws.Range("C2") "`$1,009,094,930"
Set-TextValue $ws.Range("D2") "9.60%"
Set-TextValue $ws.Range("E2") "-11.45%"

# ---- Sheet: County ----
$ws = $wb.Worksheets.Item("County")
$ws.Range("A1").Value = "Geography"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"

Set-TextValue $ws.Range("A2") "United States"
Set-TextValue $ws.Range("B2") "67.35%"
Set-TextValue $ws.Range("C2") "103,475"
Set-TextValue $ws.Range("D2") "`$267,700,640,005"
Set-TextValue $ws.Range("E2") "9.05%"
Set-TextValue $ws.Range("F2") "-12.83%"

Set-TextValue $ws.Range("A3") "Rhode Island"
Set-TextValue $ws.Range("B3") "66.00%"
Set-TextValue $ws.Range("C3") "547"
Set-TextValue $ws.Range("D3") "`$1,009,094,930"
Set-TextValue $ws.Range("E3") "9.60%"
Set-TextValue $ws.Range("F3") "-11.45%"

Set-TextValue $ws.Range("A4") "Bristol County"
Set-TextValue $ws.Range("B4") "43.48%"
Set-TextValue $ws.Range("C4") "23"
Set-TextValue $ws.Range("D4") "`$21,640,086"
Set-TextValue $ws.Range("E4") "11.00%"
Set-TextValue $ws.Range("F4") "0.89%"

Set-TextValue $ws.Range("A5") "Kent County"
Set-TextValue $ws.Range("B5") "81.48%"
Set-TextValue $ws.Range("C5") "54"
Set-TextValue $ws.Range("D5") "`$85,404,142"
Set-TextValue $ws.Range("E5") "5.56%"
Set-TextValue $ws.Range("F5") "-28.77%"

Set-TextValue $ws.Range("A6") "Newport County"
Set-TextValue $ws.Range("B6") "47.46%"
Set-TextValue $ws.Range("C6") "59"
Set-TextValue $ws.Range("D6") "`$63,188,348"
Set-TextValue $ws.Range("E6") "15.77%"
Set-TextValue $ws.Range("F6") "1.34%"

Set-TextValue $ws.Range("A7") "Providence County"
Set-TextValue $ws.Range("B7") "69.39%"
Set-TextValue $ws.Range("C7") "343"
Set-TextValue $ws.Range("D7") "`$807,901,665"
Set-TextValue $ws.Range("E7") "8.52%"
Set-TextValue $ws.Range("F7") "-13.70%"

Set-TextValue $ws.Range("A8") "Washington County"
Set-TextValue $ws.Range("B8") "60.29%"
Set-TextValue $ws.Range("C8") "68"
Set-TextValue $ws.Range("D8") "`$30,960,689"
Set-TextValue $ws.Range("E8") "14.38%"
Set-TextValue $ws.Range("F8") "-5.99%"

# ---- Sheet: Congressional District ----
$ws = $wb.Worksheets.Item("Congressional District")
$ws.Range("A1").Value = "Geography"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"

Set-TextValue $ws.Range("A2") "United States"
Set-TextValue $ws.Range("B2") "67.35%"
Set-TextValue $ws.Range("C2") "103,475"
Set-TextValue $ws.Range("D2") "`$267,700,640,005"
Set-TextValue $ws.Range("E2") "9.05%"
Set-TextValue $ws.Range("F2") "-12.83%"

Set-TextValue $ws.Range("A3") "Rhode Island"
Set-TextValue $ws.Range("B3") "66.00%"
Set-TextValue $ws.Range("C3") "547"
Set-TextValue $ws.Range("D3") "`$1,009,094,930"
Set-TextValue $ws.Range("E3") "9.60%"
Set-TextValue $ws.Range("F3") "-11.45%"

Set-TextValue $ws.Range("A4") "Congressional District 1"
Set-TextValue $ws.Range("B4") "62.28%"
Set-TextValue $ws.Range("C4") "281"
Set-TextValue $ws.Range("D4") "`$590,475,638"
Set-TextValue $ws.Range("E4") "11.27%"
Set-TextValue $ws.Range("F4") "-9.36%"

Set-TextValue $ws.Range("A5") "Congressional District 2"
Set-TextValue $ws.Range("B5") "69.92%"
Set-TextValue $ws.Range("C5") "266"
Set-TextValue $ws.Range("D5") "`$418,619,292"
Set-TextValue $ws.Range("E5") "7.96%"
Set-TextValue $ws.Range("F5") "-13.27%"

# ---- Sheet: Size ----
$ws = $wb.Worksheets.Item("Size")
$ws.Range("A1").Value = "Size"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"

Set-TextValue $ws.Range("A2") "Between `$100K and `$499K"
Set-TextValue $ws.Range("B2") "64.41%"
Set-TextValue $ws.Range("C2") "177"
Set-TextValue $ws.Range("D2") "`$20,163,196"
Set-TextValue $ws.Range("E2") "14.24%"
Set-TextValue $ws.Range("F2") "-12.83%"

Set-TextValue $ws.Range("A3") "Between `$1M and `$4.99M"
Set-TextValue $ws.Range("B3") "68.39%"
Set-TextValue $ws.Range("C3") "155"
Set-TextValue $ws.Range("D3") "`$150,835,594"
Set-TextValue $ws.Range("E3") "11.45%"
Set-TextValue $ws.Range("F3") "-17.52%"

Set-TextValue $ws.Range("A4") "Between `$500K and `$999K"
Set-TextValue $ws.Range("B4") "65.82%"
Set-TextValue $ws.Range("C4") "79"
Set-TextValue $ws.Range("D4") "`$17,881,469"
Set-TextValue $ws.Range("E4") "8.94%"
Set-TextValue $ws.Range("F4") "-13.60%"

Set-TextValue $ws.Range("A5") "Between `$5M and `$9.99M"
Set-TextValue $ws.Range("B5") "71.11%"
Set-TextValue $ws.Range("C5") "45"
Set-TextValue $ws.Range("D5") "`$92,176,427"
Set-TextValue $ws.Range("E5") "7.68%"
Set-TextValue $ws.Range("F5") "-8.25%"

Set-TextValue $ws.Range("A6") "Greater than `$10M"
Set-TextValue $ws.Range("B6") "70.31%"
Set-TextValue $ws.Range("C6") "64"
Set-TextValue $ws.Range("D6") "`$726,810,103"
Set-TextValue $ws.Range("E6") "4.87%"
Set-TextValue $ws.Range("F6") "-6.99%"

Set-TextValue $ws.Range("A7") "Less than `$100K"
Set-TextValue $ws.Range("B7") "44.44%"
Set-TextValue $ws.Range("C7") "27"
Set-TextValue $ws.Range("D7") "`$1,228,141"
Set-TextValue $ws.Range("E7") "24.29%"
Set-TextValue $ws.Range("F7") "3.70%"

Set-TextValue $ws.Range("A8") "Total"
Set-TextValue $ws.Range("B8") "66.00%"
Set-TextValue $ws.Range("C8") "547"
Set-TextValue $ws.Range("D8") "`$1,009,094,930"
Set-TextValue $ws.Range("E8") "9.60%"
Set-TextValue $ws.Range("F8") "-11.45%"

# ---- Sheet: Subsector ----
$ws = $wb.Worksheets.Item("Subsector")
$ws.Range("A1").Value = "Subsector"
$ws.Range("B1").Value = "Share of 990 filers with government grants at risk"
$ws.Range("C1").Value = "Number of 990 filers with government grants"
$ws.Range("D1").Value = "Total government grants (`$)"
$ws.Range("E1").Value = "Size of operating surplus with government grants"
$ws.Range("F1").Value = "Size of operating surplus without government grants"

Set-TextValue $ws.Range("A2") "Arts, Culture, and Humanities"
Set-TextValue $ws.Range("B2") "58.73%"
Set-TextValue $ws.Range("C2") "63"
Set-TextValue $ws.Range("D2") "`$28,568,314"
Set-TextValue $ws.Range("E2") "19.67%"
Set-TextValue $ws.Range("F2") "-2.35%"

Set-TextValue $ws.Range("A3") "Education (Excluding Universities)"
Set-TextValue $ws.Range("B3") "68.29%"
Set-TextValue $ws.Range("C3") "82"
Set-TextValue $ws.Range("D3") "`$74,378,415"
Set-TextValue $ws.Range("E3") "7.19%"
Set-TextValue $ws.Range("F3") "-12.89%"

Set-TextValue $ws.Range("A4") "Environment and Animals"
Set-TextValue $ws.Range("B4") "34.62%"
Set-TextValue $ws.Range("C4") "26"
Set-TextValue $ws.Range("D4") "`$6,959,923"
Set-TextValue $ws.Range("E4") "30.08%"
Set-TextValue $ws.Range("F4") "13.73%"

Set-TextValue $ws.Range("A5") "Health (Excluding Hospitals)"
Set-TextValue $ws.Range("B5") "78.57%"
Set-TextValue $ws.Range("C5") "42"
Set-TextValue $ws.Range("D5") "`$77,536,872"
Set-TextValue $ws.Range("E5") "5.19%"
Set-TextValue $ws.Range("F5") "-17.02%"

Set-TextValue $ws.Range("A6") "Hospitals"
Set-TextValue $ws.Range("B6") "100.00%"
Set-TextValue $ws.Range("C6") "2"
Set-TextValue $ws.Range("D6") "`$5,765,668"
Set-TextValue $ws.Range("E6") "5.24%"
Set-TextValue $ws.Range("F6") "-35.54%"

Set-TextValue $ws.Range("A7") "Human Services"
Set-TextValue $ws.Range("B7") "65.36%"
Set-TextValue $ws.Range("C7") "153"
Set-TextValue $ws.Range("D7") "`$158,026,958"
Set-TextValue $ws.Range("E7") "7.86%"
Set-TextValue $ws.Range("F7") "-16.79%"

Set-TextValue $ws.Range("A8") "International, Foreign Affairs"
Set-TextValue $ws.Range("B8") "60.00%"
Set-TextValue $ws.Range("C8") "5"
Set-TextValue $ws.Range("D8") "`$19,732,448"
Set-TextValue $ws.Range("E8") "2.05%"
Set-TextValue $ws.Range("F8") "-15.04%"

Set-TextValue $ws.Range("A9") "Mutual/Membership Benefit"
Set-TextValue $ws.Range("B9") "100.00%"
Set-TextValue $ws.Range("C9") "2"
Set-TextValue $ws.Range("D9") "`$192,493"
Set-TextValue $ws.Range("E9") "1.68%"
Set-TextValue $ws.Range("F9") "-11.58%"

Set-TextValue $ws.Range("A10") "Public, Societal Benefit"
Set-TextValue $ws.Range("B10") "73.53%"
Set-TextValue $ws.Range("C10") "34"
Set-TextValue $ws.Range("D10") "`$37,515,527"
Set-TextValue $ws.Range("E10") "5.20%"
Set-TextValue $ws.Range("F10") "-38.05%"

Set-TextValue $ws.Range("A11") "Religion Related"
Set-TextValue $ws.Range("B11") "66.67%"
Set-TextValue $ws.Range("C11") "3"
Set-TextValue $ws.Range("D11") "`$436,092"
Set-TextValue $ws.Range("E11") "12.47%"
Set-TextValue $ws.Range("F11") "-11.42%"

Set-TextValue $ws.Range("A12") "Unclassified"
Set-TextValue $ws.Range("B12") "69.77%"
Set-TextValue $ws.Range("C12") "129"
Set-TextValue $ws.Range("D12") "`$502,581,531"
Set-TextValue $ws.Range("E12") "9.36%"
Set-TextValue $ws.Range("F12") "-13.61%"

Set-TextValue $ws.Range("A13") "Universities"
Set-TextValue $ws.Range("B13") "33.33%"
Set-TextValue $ws.Range("C13") "6"
Set-TextValue $ws.Range("D13") "`$97,400,689"
Set-TextValue $ws.Range("E13") "7.53%"
Set-TextValue $ws.Range("F13") "2.68%"

Set-TextValue $ws.Range("A14") "Total"
Set-TextValue $ws.Range("B14") "66.00%"
Set-TextValue $ws.Range("C14") "547"
Set-TextValue $ws.Range("D14") "`$1,009,094,930"
Set-TextValue $ws.Range("E14") "9.60%"
Set-TextValue $ws.Range("F14") "-11.45%"

